$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13's politeness_score (B13) becomes a genuine number (it previously held "3" as text)
$ws.Cells.Item(13, 2).Value = 3

# New row 14: another annotation entry for Ruilin, appended after row 13
$ws.Cells.Item(14, 1).Value = "Ruilin"

# Keep politeness_score as text "3" (matches source data), not auto-converted to a number
$ws.Cells.Item(14, 2).NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = "3"
$ws.Cells.Item(14, 2).Style = "Normal"

$ws.Cells.Item(14, 3).Value = "无"
$ws.Cells.Item(14, 4).Value = "SUG"
$ws.Cells.Item(14, 5).Value = "WRI"
$ws.Cells.Item(14, 6).Value = "3a6bf25f-9f71-48b7-a40b-7e968e5f9337"
$ws.Cells.Item(14, 7).Value = "ry-TW-WAb_annotated.xlsx"
$ws.Cells.Item(14, 8).Value = "I suggest to change it to e.g. 'from the true to the approximate posterior' to avoid confusion."
